$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 72

$ws.Cells.Item($row, 1).Value = "2025-11-27"
$ws.Cells.Item($row, 2).Value = "Pick 3"
$ws.Cells.Item($row, 3).Value = "251127"
$ws.Cells.Item($row, 4).Value = "6-7-6"
$ws.Cells.Item($row, 5).Value = "2025-11-27T21:37:50.045+04:00"
